$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.56"
$ws.Range("E2").Value = "'1.10%"
$ws.Range("D3").Value = "'31.89"
$ws.Range("E3").Value = "'1.76%"
$ws.Range("D4").Value = "'5.119"
$ws.Range("E4").Value = "'0.28%"
$ws.Range("D5").Value = "'0.07841"
$ws.Range("E5").Value = "'-2.32%"
$ws.Range("D6").Value = "'2.246"
$ws.Range("E6").Value = "'-8.08%"
$ws.Range("D7").Value = "'7.811"
$ws.Range("E7").Value = "'-0.13%"
$ws.Range("D8").Value = "'3.810"
$ws.Range("E8").Value = "'0.33%"
$ws.Range("D9").Value = "'0.9279"
$ws.Range("E9").Value = "'0.84%"
$ws.Range("D10").Value = "'0.1771"
$ws.Range("E10").Value = "'2.26%"
$ws.Range("D11").Value = "'0.07653"
$ws.Range("E11").Value = "'4.20%"
$ws.Range("D12").Value = "'0.08846"
$ws.Range("E12").Value = "'2.86%"
$ws.Range("E13").Value = "'2.14%"
$ws.Range("E14").Value = "'0.52%"
$ws.Range("D15").Value = "'0.001517"
$ws.Range("E15").Value = "'0.73%"
$ws.Range("D16").Value = "'0.005865"
$ws.Range("E16").Value = "'-2.51%"
$ws.Range("D17").Value = "'3.463"
$ws.Range("E17").Value = "'-1.37%"
$ws.Range("D18").Value = "'2.251"
$ws.Range("E18").Value = "'0.16%"
$ws.Range("D21").Value = "'4.324"
$ws.Range("E21").Value = "'-6.50%"
$ws.Range("E22").Value = "'10.69%"
$ws.Range("D23").Value = "'0.04603"
$ws.Range("E23").Value = "'-0.62%"
$ws.Range("E24").Value = "'0.36%"
$ws.Range("D25").Value = "'0.004484"
$ws.Range("E25").Value = "'1.17%"
$ws.Range("D26").Value = "'0.0001250"
$ws.Range("E26").Value = "'4.05%"
$ws.Range("E27").Value = "'-1.40%"
$ws.Range("D39").Value = "'0.01782"
$ws.Range("E39").Value = "'-0.73%"
$ws.Range("D40").Value = "'0.04781"
$ws.Range("E40").Value = "'6.01%"
$ws.Range("D41").Value = "'0.007358"
$ws.Range("E41").Value = "'4.80%"
$ws.Range("D42").Value = "'0.1363"
$ws.Range("E42").Value = "'1.58%"
$ws.Range("D43").Value = "'0.002190"
$ws.Range("E43").Value = "'-2.35%"
$ws.Range("D44").Value = "'0.009847"
$ws.Range("E44").Value = "'-0.03%"
$ws.Range("E45").Value = "'-5.05%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.12%"
$ws.Range("D48").Value = "'0.7006"
$ws.Range("E48").Value = "'-14.62%"
$ws.Range("E49").Value = "'-0.12%"
$ws.Range("E50").Value = "'-0.12%"
